$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data table shifted up by one row (the oldest observation dropped off
# the front of the rolling window) and the forecasts were recomputed.
# Remove the now-unused last row (row 19) - this also updates the sheet
# dimension from A1:E19 to A1:E18.
$ws.Rows(19).Delete()

# Column A: date_of_forecast (style s=2, already present on the cells)
$ws.Range("A2").Value = 39765
$ws.Range("A3").Value = 40130
$ws.Range("A4").Value = 40494
$ws.Range("A5").Value = 40862
$ws.Range("A6").Value = 41228
$ws.Range("A7").Value = 41592
$ws.Range("A8").Value = 41957
$ws.Range("A9").Value = 42321
$ws.Range("A10").Value = 42689
$ws.Range("A11").Value = 43053
$ws.Range("A12").Value = 43418
$ws.Range("A13").Value = 43783
$ws.Range("A14").Value = 44159
$ws.Range("A15").Value = 44525
$ws.Range("A16").Value = 44890
$ws.Range("A17").Value = 45254
$ws.Range("A18").Value = 45618

# Column B: y_0
$ws.Range("B2").Value = 2008
$ws.Range("B3").Value = 2009
$ws.Range("B4").Value = 2010
$ws.Range("B5").Value = 2011
$ws.Range("B6").Value = 2012
$ws.Range("B7").Value = 2013
$ws.Range("B8").Value = 2014
$ws.Range("B9").Value = 2015
$ws.Range("B10").Value = 2016
$ws.Range("B11").Value = 2017
$ws.Range("B12").Value = 2018
$ws.Range("B13").Value = 2019
$ws.Range("B14").Value = 2020
$ws.Range("B15").Value = 2021
$ws.Range("B16").Value = 2022
$ws.Range("B17").Value = 2023
$ws.Range("B18").Value = 2024

# Column C: y_0_forecast (row 2 has no forecast, same as before)
$ws.Range("C3").Value = -1.324983933426893
$ws.Range("C4").Value = -0.3900454704678369
$ws.Range("C5").Value = -0.29958481534893
$ws.Range("C6").Value = -0.2075757021743008
$ws.Range("C7").Value = 0.124712275190686
$ws.Range("C8").Value = -0.255298189276465
$ws.Range("C9").Value = 0.07418514192796266
$ws.Range("C10").Value = -0.07611406013281474
$ws.Range("C11").Value = -0.191300579729714
$ws.Range("C12").Value = 0.0970330232288763
$ws.Range("C13").Value = -0.7407518902333265
$ws.Range("C14").Value = 0.3056679541520335
$ws.Range("C15").Value = -1.388491535160907
$ws.Range("C16").Value = -1.678482969789596
$ws.Range("C17").Value = -0.5999457276250508
$ws.Range("C18").Value = -0.05499271238530445

# Column D: y_1
$ws.Range("D2").Value = 2009
$ws.Range("D3").Value = 2010
$ws.Range("D4").Value = 2011
$ws.Range("D5").Value = 2012
$ws.Range("D6").Value = 2013
$ws.Range("D7").Value = 2014
$ws.Range("D8").Value = 2015
$ws.Range("D9").Value = 2016
$ws.Range("D10").Value = 2017
$ws.Range("D11").Value = 2018
$ws.Range("D12").Value = 2019
$ws.Range("D13").Value = 2020
$ws.Range("D14").Value = 2021
$ws.Range("D15").Value = 2022
$ws.Range("D16").Value = 2023
$ws.Range("D17").Value = 2024
$ws.Range("D18").Value = 2025

# Column E: y_1_forecast. The AR(2) model now needs two prior observations
# of y_1 before it can produce a forecast, so E2:E5 no longer hold values
# (previously E3:E5 were populated); clear them and fill in the
# recomputed forecasts for rows 6-18.
$ws.Range("E2:E5").ClearContents()
$ws.Range("E6").Value = -0.4223781730902543
$ws.Range("E7").Value = -0.3577371449824729
$ws.Range("E8").Value = -0.2877346565283379
$ws.Range("E9").Value = -0.2097319935285391
$ws.Range("E10").Value = -0.3036570471216304
$ws.Range("E11").Value = -0.1907914728172644
$ws.Range("E12").Value = -0.2529765062333933
$ws.Range("E13").Value = -0.428077259747528
$ws.Range("E14").Value = -0.2112001730687485
$ws.Range("E15").Value = -1.701252732314051
$ws.Range("E16").Value = -0.6241481568271312
$ws.Range("E17").Value = -0.2119687890143274
$ws.Range("E18").Value = -0.06564014165270082
